# Jenkins automation run: refresh the WebData query text in Sheet1!B5
# (SQL reformatted/cleaned up, file-size rounding switched from ROUND to a
# trimmed printf, "Left JOIN" -> "LEFT JOIN", and an ORDER BY / LIMIT added),
# then leave the active selection on the cell that was edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$sql = @"
WITH file_data AS (
    SELECT file_name, data_category, file_type, file_size, file_access, file_description, "sample.id" 
    FROM df_sequencing_file
    UNION
    SELECT file_name, data_category, file_type, file_size, file_access, file_description, "sample.id" 
    FROM df_pathology_file
)
SELECT DISTINCT
    fd.file_name AS "File Name",
    fd.data_category AS "Data Category",
    COALESCE(fd.file_description, '') AS "File Description",
    fd.file_type AS "File Type",
    CASE     
        WHEN fd.file_size >= 1024 * 1024 * 1024 THEN 
            RTRIM(RTRIM(printf('%.2f', fd.file_size / (1024.0 * 1024.0 * 1024.0)), '0'), '.') || ' GB'
        WHEN fd.file_size >= 1024 * 1024 THEN 
            RTRIM(RTRIM(printf('%.2f', fd.file_size / (1024.0 * 1024.0)), '0'), '.') || ' MB'
        WHEN fd.file_size >= 1024 THEN 
            RTRIM(RTRIM(printf('%.2f', fd.file_size / 1024.0), '0'), '.') || ' KB'
        ELSE 
            RTRIM(RTRIM(printf('%.2f', fd.file_size), '0'), '.') || ' Bytes'
    END AS "File Size",
    fd.file_access AS "File Access",
    std.dbgap_accession AS "Study ID",
    prt.participant_id AS "Participant ID",
    smp.sample_id AS "Sample ID"    
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_sample smp ON prt.id = smp."participant.id"
JOIN 
    file_data fd ON smp.id = fd."sample.id"
LEFT JOIN 
    df_diagnosis diag ON prt.id = diag."participant.id"     
LEFT JOIN 
    df_survival srv ON prt.id = srv."participant.id"
WHERE 
    std.dbgap_accession = 'phs000720' 
    AND diag.diagnosis = '8900/3 : Rhabdomyosarcoma, NOS' 
    AND srv.last_known_survival_status = 'Not Reported'
ORDER BY 
    fd.file_name ASC 
LIMIT 100;
"@

$ws.Range("B5").Value = $sql

# Match the author's final cursor position (was A11, now B5 - the cell that
# was just edited).
[void]$ws.Range("B5").Select()
